# Generate Report for Handback
# - Overview/zh-cn/de-de "Status" cells for the 1690f5f7... row flip from
#   "Ready for handoff" to "Handback transform failed".
# - The zh-cn and de-de "Error Detail" cells (column P) for that row get a
#   handback-transform error message, and column P is widened to fit it.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C3").Value = $newStatus
$ws2.Range("P3").Value = "Handback file name: nodilwnx.zyj is different with handoff file name: 1690f5f7-6bd0-46ed-84fb-e07da543fa41.5717fd8150bfbfdb6c152c665b6d75afe1ec070d.zh-cn."
$ws2.Columns.Item(16).ColumnWidth = 39.17

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C3").Value = $newStatus
$ws3.Range("P3").Value = "Handback file name: nodilwnx.zyj is different with handoff file name: 1690f5f7-6bd0-46ed-84fb-e07da543fa41.5717fd8150bfbfdb6c152c665b6d75afe1ec070d.de-de."
$ws3.Columns.Item(16).ColumnWidth = 39.17
